# Update the division problems in the table to the new values.
# Each source string is unique within the document, and replacements are
# applied in the same order they appear in the document so that a newly
# created value (e.g. "23÷5=" produced later) is never re-matched by an
# earlier rule that already ran.

$d = $word.ActiveDocument

$replacements = @(
    @("23÷6=", "12÷5="),
    @("23÷5=", "40÷9="),
    @("63÷5=", "36÷3="),
    @("24÷5=", "24÷8="),
    @("53÷5=", "79÷7="),
    @("32÷4=", "23÷5="),
    @("73÷5=", "63÷3="),
    @("62÷7=", "11÷5="),
    @("73÷2=", "94÷3="),
    @("91÷5=", "33÷3="),
    @("93÷5=", "53÷7="),
    @("18÷5=", "73÷5="),
    @("29÷4=", "68÷3="),
    @("58÷7=", "32÷7="),
    @("37÷3=", "71÷5="),
    @("36÷9=", "55÷4="),
    @("73÷3=", "81÷7="),
    @("96÷7=", "98÷4="),
    @("23÷2=", "28÷3="),
    @("43÷9=", "42÷5="),
    @("31÷3=", "78÷3="),
    @("92÷6=", "20÷4="),
    @("74÷2=", "13÷3="),
    @("62÷8=", "19÷7="),
    @("59÷2=", "72÷2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $new, 2)
}

$d.Save()
